$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# True-up of past four days of discharges (column E) that were previously blank.
$ws.Range("E24").Value2 = 1763
$ws.Range("E25").Value2 = 1944
$ws.Range("E26").Value2 = 1952

# Add New York state hospitalization data for 10 and 11 April 2020 (rows 27-28).
# Copy formatting from the last existing row first so the new rows inherit the
# same cell styles (e.g. date formatting in column A).
$ws.Range("A26:F26").Copy()
$ws.Range("A27:F27").PasteSpecial(-4122)
$ws.Range("A27:F27").Copy()
$ws.Range("A28:F28").PasteSpecial(-4122)

# Row 27: 10 April 2020
$ws.Range("A27").Value2 = 43931
$ws.Range("B27").Value2 = 85
$ws.Range("C27").Value2 = 101
$ws.Range("D27").Value2 = -26
$ws.Range("E27").Value2 = 1776
$ws.Range("F27").Value2 = 783

# Row 28: 11 April 2020
$ws.Range("A28").Value2 = 43932
$ws.Range("B28").Value2 = 53
$ws.Range("C28").Value2 = 189
$ws.Range("D28").Value2 = 110
$ws.Range("E28").Value2 = 1862
$ws.Range("F28").Value2 = 758

# Match the author's final selection/active cell.
$ws.Range("F29").Select()
